$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2255.7144
$ws.Range("I62").Value = 1960
$ws.Range("J62").Value = 2995
$ws.Range("K62").Value = 1960
$ws.Range("L62").Value = 2995
$ws.Range("M62").Value = -1336
$ws.Range("N62").Value = -4243
$ws.Range("H65").Value = 2255.7144
$ws.Range("I65").Value = 1960
$ws.Range("J65").Value = 2995
$ws.Range("K65").Value = 9800
$ws.Range("L65").Value = 14975
$ws.Range("M65").Value = -6680
$ws.Range("N65").Value = -21215
$ws.Range("H132").Value = 4036056
$ws.Range("I132").Value = 4549411
$ws.Range("J132").Value = 2550.7144
$ws.Range("K132").Value = 13648233
$ws.Range("L132").Value = 7652.1432
$ws.Range("M132").Value = -13645703
$ws.Range("N132").Value = -12712.1432
$ws.Range("H137").Value = 1600.6285
$ws.Range("I137").Value = 1130.0714
$ws.Range("J137").Value = 3482.8572
$ws.Range("K137").Value = 3390.2142
$ws.Range("L137").Value = 10448.5716
$ws.Range("M137").Value = -840.2142000000003
$ws.Range("N137").Value = -15548.5716

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 36489.168
$ws.Range("I45").Value = 46278.22
$ws.Range("K45").Value = 46278.22
$ws.Range("M45").Value = -45901.22
$ws.Range("H63").Value = 2733.3333
$ws.Range("I63").Value = 1850
$ws.Range("J63").Value = 3616.6667
$ws.Range("K63").Value = 1850
$ws.Range("L63").Value = 3616.6667
$ws.Range("M63").Value = -1164
$ws.Range("N63").Value = -4988.6667
$ws.Range("H66").Value = 2733.3333
$ws.Range("I66").Value = 1850
$ws.Range("J66").Value = 3616.6667
$ws.Range("K66").Value = 9250
$ws.Range("L66").Value = 18083.3335
$ws.Range("M66").Value = -5818
$ws.Range("N66").Value = -24947.3335

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 2140.84
$ws.Range("I80").Value = 884.125
$ws.Range("J80").Value = 2732.2354
$ws.Range("K80").Value = 884.125
$ws.Range("L80").Value = 2732.2354
$ws.Range("M80").Value = 113.875
$ws.Range("N80").Value = -4728.2354
$ws.Range("H83").Value = 2140.84
$ws.Range("I83").Value = 884.125
$ws.Range("J83").Value = 2732.2354
$ws.Range("K83").Value = 4420.625
$ws.Range("L83").Value = 13661.177
$ws.Range("M83").Value = 571.375
$ws.Range("N83").Value = -23645.177
$ws.Range("H119").Value = 42974.5
$ws.Range("J119").Value = 42974.5
$ws.Range("L119").Value = 42974.5
$ws.Range("N119").Value = -52650.5
$ws.Range("H132").Value = 67352.82000000001
$ws.Range("J132").Value = 67352.82000000001
$ws.Range("L132").Value = 67352.82000000001
$ws.Range("N132").Value = -77472.82000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 3100
$ws.Range("I4").Value = 1250
$ws.Range("J4").Value = 4333.3335
$ws.Range("K4").Value = 1250
$ws.Range("L4").Value = 4333.3335
$ws.Range("M4").Value = -1138
$ws.Range("N4").Value = -4557.3335
$ws.Range("H16").Value = 1132.2
$ws.Range("I16").Value = 1166.5
$ws.Range("K16").Value = 1166.5
$ws.Range("M16").Value = -879.5
$ws.Range("H31").Value = 27574.344
$ws.Range("I31").Value = 1234.6154
$ws.Range("J31").Value = 47141
$ws.Range("K31").Value = 1234.6154
$ws.Range("L31").Value = 47141
$ws.Range("M31").Value = -939.6153999999999
$ws.Range("N31").Value = -47731
$ws.Range("H34").Value = 27574.344
$ws.Range("I34").Value = 1234.6154
$ws.Range("J34").Value = 47141
$ws.Range("K34").Value = 1234.6154
$ws.Range("L34").Value = 47141
$ws.Range("M34").Value = -1032.6154
$ws.Range("N34").Value = -47545
$ws.Range("H86").Value = 2141.5
$ws.Range("I86").Value = 1692.1428
$ws.Range("J86").Value = 2770.6
$ws.Range("K86").Value = 1692.1428
$ws.Range("L86").Value = 2770.6
$ws.Range("M86").Value = -569.1428000000001
$ws.Range("N86").Value = -5016.6
$ws.Range("H89").Value = 2141.5
$ws.Range("I89").Value = 1692.1428
$ws.Range("J89").Value = 2770.6
$ws.Range("K89").Value = 8460.714
$ws.Range("L89").Value = 13853
$ws.Range("M89").Value = -2844.714
$ws.Range("N89").Value = -25085
$ws.Range("H105").Value = 2169.5
$ws.Range("I105").Value = 2188.6
$ws.Range("J105").Value = 2074
$ws.Range("K105").Value = 2188.6
$ws.Range("L105").Value = 2074
$ws.Range("M105").Value = -441.5999999999999
$ws.Range("N105").Value = -5568
$ws.Range("H113").Value = 1132.2
$ws.Range("I113").Value = 1166.5
$ws.Range("K113").Value = 1166.5
$ws.Range("M113").Value = 1003.5
$ws.Range("H132").Value = 3807.513
$ws.Range("I132").Value = 3787.48
$ws.Range("J132").Value = 3843.2856
$ws.Range("K132").Value = 11362.44
$ws.Range("L132").Value = 11529.8568
$ws.Range("M132").Value = -8832.440000000001
$ws.Range("N132").Value = -16589.8568

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 47266.426
$ws.Range("I70").Value = 72969.24000000001
$ws.Range("J70").Value = 5856.3335
$ws.Range("K70").Value = 72969.24000000001
$ws.Range("L70").Value = 5856.3335
$ws.Range("M70").Value = -72699.24000000001
$ws.Range("N70").Value = -6396.3335
$ws.Range("H73").Value = 47266.426
$ws.Range("I73").Value = 72969.24000000001
$ws.Range("J73").Value = 5856.3335
$ws.Range("K73").Value = 72969.24000000001
$ws.Range("L73").Value = 5856.3335
$ws.Range("M73").Value = -72033.24000000001
$ws.Range("N73").Value = -7728.3335
$ws.Range("H97").Value = 66669452
$ws.Range("I97").Value = 76925840
$ws.Range("J97").Value = 2925
$ws.Range("K97").Value = 76925840
$ws.Range("L97").Value = 2925
$ws.Range("M97").Value = -76925344
$ws.Range("N97").Value = -3917
$ws.Range("H107").Value = 421179.72
$ws.Range("I107").Value = 308.6
$ws.Range("J107").Value = 2525535.2
$ws.Range("K107").Value = 308.6
$ws.Range("L107").Value = 2525535.2
$ws.Range("M107").Value = 1611.4
$ws.Range("N107").Value = -2529375.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 83271.92999999999
$ws.Range("I2").Value = 375375
$ws.Range("J2").Value = 5377.778
$ws.Range("K2").Value = 375375
$ws.Range("L2").Value = 5377.778
$ws.Range("M2").Value = -375263
$ws.Range("N2").Value = -5601.778
$ws.Range("H68").Value = 5998
$ws.Range("I68").Value = 4334
$ws.Range("J68").Value = 6830
$ws.Range("K68").Value = 4334
$ws.Range("L68").Value = 6830
$ws.Range("M68").Value = -3585
$ws.Range("N68").Value = -8328
$ws.Range("H71").Value = 5998
$ws.Range("I71").Value = 4334
$ws.Range("J71").Value = 6830
$ws.Range("K71").Value = 21670
$ws.Range("L71").Value = 34150
$ws.Range("M71").Value = -17926
$ws.Range("N71").Value = -41638

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 4170000
$ws.Range("J5").Value = 6250000
$ws.Range("L5").Value = 6250000
$ws.Range("N5").Value = -6250224
$ws.Range("H107").Value = 71835.64
$ws.Range("I107").Value = 335
$ws.Range("J107").Value = 250587.25
$ws.Range("K107").Value = 1005
$ws.Range("L107").Value = 751761.75
$ws.Range("M107").Value = 915
$ws.Range("N107").Value = -755601.75
$ws.Range("H113").Value = 870.2
$ws.Range("I113").Value = 633.3333
$ws.Range("J113").Value = 1225.5
$ws.Range("K113").Value = 1899.9999
$ws.Range("L113").Value = 3676.5
$ws.Range("M113").Value = 270.0001
$ws.Range("N113").Value = -8016.5
$ws.Range("H123").Value = 29847.5
$ws.Range("J123").Value = 29847.5
$ws.Range("L123").Value = 29847.5
$ws.Range("N123").Value = -39647.5
